# Changing test case name for watch list scripts
# Renames the TCID values (column A, rows 2-30) on the "Test Cases" sheet
# from TestCase_EN to sequential Watchlist0NN identifiers, and leaves the
# final selection on D30 (matching the author's last-saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Select()

$names = @(
    "Watchlist001","Watchlist002","Watchlist003","Watchlist004","Watchlist005",
    "Watchlist006","Watchlist007","Watchlist008","Watchlist009","Watchlist010",
    "Watchlist011","Watchlist012","Watchlist013","Watchlist014","Watchlist015",
    "Watchlist016","Watchlist017","Watchlist018","Watchlist019","Watchlist020",
    "Watchlist021","Watchlist022","Watchlist023","Watchlist024","Watchlist025",
    "Watchlist026","Watchlist027","Watchlist028","Watchlist029"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $names[$i]
}

# Match the saved cursor/scroll state from the commit: active cell D30.
$ws.Range("D30").Select()
